$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.996.15'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '3.155.75'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'573.98"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.22%  '
$ws.Range('D6').Value = "'149.97"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.60%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '3.152.26'
$ws.Range('E8').Value = '  +2.77%  '
$ws.Range('E9').Value = '  +1.84%  '
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').Value = "'0.498"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.01%  '
$ws.Range('D13').Value = "'0.0000266"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +14.15%  '
$ws.Range('D14').Value = "'37.14"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.05%  '
$ws.Range('D15').Value = '3.674.39'
$ws.Range('E15').Value = '  +3.06%  '
$ws.Range('D16').Value = '65.047.68'
$ws.Range('E16').Value = '  +1.65%  '
$ws.Range('D17').Value = '3.154.30'
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').Value = "'7.10"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.77%  '
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').Value = "'506.09"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.51%  '
$ws.Range('D21').Value = "'14.85"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.20%  '
$ws.Range('D22').Value = "'0.718"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.23%  '
$ws.Range('D23').Value = "'15.34"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.03%  '
$ws.Range('E24').Value = '  +2.25%  '
$ws.Range('D25').Value = "'84.33"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.02%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = "'2.91"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.96%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = "'8.90"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.85%  '
$ws.Range('D29').Value = "'2.17"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.17%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').Value = "'27.61"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.41%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').Value = "'2.78"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.64%  '
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('E33').Value = '  +2.48%  '
$ws.Range('D34').Value = "'6.17"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.70%  '
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('D36').Value = "'54.77"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('D37').Value = "'0.0898"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.28%  '
$ws.Range('D38').Value = "'464.98"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.38%  '
$ws.Range('D39').Value = "'0.0420"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').Value = "'2.98"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.83%  '
$ws.Range('E41').Value = '  +3.39%  '
$ws.Range('D42').Value = '3.049.74'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = "'0.117"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('E44').Value = '  +8.30%  '
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('D46').Value = "'28.52"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.25%  '
$ws.Range('D47').Value = '0.0₃0587'
$ws.Range('E47').Value = '  +11.76%  '
$ws.Range('D49').Value = "'0.114"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('E50').Value = '  +4.52%  '
$ws.Range('D51').Value = "'119.43"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.72%  '
